$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: names and card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long, purely-numeric-looking card number that must stay as
# TEXT (as in the original file) rather than being auto-converted to a
# number. Force text entry with a leading apostrophe, then re-apply the
# original direct formatting (from the untouched D3 cell, which shares
# B3's style) so no stray number formatting sticks to the cell.
$ws.Range("B3").Formula = "'2570314725427075"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null

$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 11.11.2024"

# Row 6: transaction 1
$ws.Range("B6").Value = "13.11."
$ws.Range("C6").Value = "14.11."
$ws.Range("D6").Value = "PAYPAL LDGAKP"
$ws.Range("E6").Value = "77,00-"

# Row 7: transaction 2
$ws.Range("B7").Value = "16.11."
$ws.Range("C7").Value = "17.11."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "25,05-"

# Row 8: transaction 3
$ws.Range("B8").Value = "19.11."
$ws.Range("C8").Value = "20.11."
$ws.Range("D8").Value = "ZALANDO MKTPLC EU EEPVKL"
$ws.Range("E8").Value = "53,13-"

# Row 9: new transaction 4 (previously empty)
# E9 needs to switch from the "center" style (s=13) used by blank rows
# to the "right-aligned" style (s=17) used by the other amount cells.
# Copy the exact format from E8 (same style) before setting the value.
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null

$ws.Range("B9").Value = "22.11."
$ws.Range("C9").Value = "23.11."
$ws.Range("D9").Value = "PAYPAL NJVJRV"
$ws.Range("E9").Value = "44,88-"

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 25.11.2024"
$ws.Range("E12").Value = "200,06-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 30.11.2024"
